$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B1").Value = "One2172140897"
$ws.Range("B4").Value = '<strong> Дата: 16-04-2020</strong></p>'
$ws.Range("B5").Value = '<p> <li style="list-style-type: none;" >                      <a href="#" class="" style="padding: 0px" data-toggle="dropdown" role="button"                      aria-haspopup="true" aria-expanded="false"><strong>Инфо</strong>                      <span class="caret"></span></a> <ul class="dropdown-menu">                   <li>info</li></ul> </li>'
